$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row to the "Condicion_Pacientes" Excel table with the
# latest day's figures (2020-05-13 / serial 43964).
$lo = $ws.ListObjects.Item("Condicion_Pacientes")
$newRow = $lo.ListRows.Add()
$rng = $newRow.Range

# Copy the date cell above first so the new date cell inherits the same
# date number format/style, then overwrite the values.
$ws.Range("A44").Copy($rng.Item(1))

$rng.Item(1).Value = 43964
$rng.Item(2).Value = 451
$rng.Item(3).Value = 175
$rng.Item(4).Value = 243
$rng.Item(5).Value = 7
$rng.Item(6).Value = 21

# Match the author's final selection/cursor position after entering the
# new row of data.
$ws.Range("A46").Select() | Out-Null
